# --- Update the "cryptos" market snapshot (prices + 1h volume%) ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column ("Price") cells hold text (values like "70.265.30" or
# "0.0000122" are not valid numbers), so force Text format before
# writing, then restore the default style so no formatting is left
# behind on these cells.
$dCells = @("D2","D3","D4","D5","D6","D7","D9","D10","D11","D12","D13","D15","D16","D17","D19","D20","D21","D23","D24","D25","D26","D27","D29","D30","D31","D32","D33","D35","D36","D38","D39","D41","D42","D43","D44","D45","D46","D48","D49","D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "70.362.72"
$ws.Range("E2").Value = "  +4.71%  "

# Row 3
$ws.Range("D3").Value = "3.788.80"
$ws.Range("E3").Value = "  +22.01%  "

# Row 4
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "617.95"
$ws.Range("E5").Value = "  +7.69%  "

# Row 6
$ws.Range("D6").Value = "177.97"
$ws.Range("E6").Value = "  +0.30%  "

# Row 7
$ws.Range("D7").Value = "3.786.72"
$ws.Range("E7").Value = "  +21.98%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").Value = "0.553"
$ws.Range("E9").Value = "  +6.94%  "

# Row 10
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +10.86%  "

# Row 11
$ws.Range("D11").Value = "6.40"
$ws.Range("E11").Value = "  -1.80%  "

# Row 12
$ws.Range("D12").Value = "0.505"
$ws.Range("E12").Value = "  +8.13%  "

# Row 13
$ws.Range("D13").Value = "40.82"
$ws.Range("E13").Value = "  +12.75%  "

# Row 14
$ws.Range("E14").Value = "  +7.27%  "

# Row 15
$ws.Range("D15").Value = "4.417.43"
$ws.Range("E15").Value = "  +22.04%  "

# Row 16
$ws.Range("D16").Value = "3.779.20"
$ws.Range("E16").Value = "  +22.00%  "

# Row 17
$ws.Range("D17").Value = "70.473.04"
$ws.Range("E17").Value = "  +5.03%  "

# Row 18
$ws.Range("E18").Value = "  +1.39%  "

# Row 19
$ws.Range("D19").Value = "7.64"
$ws.Range("E19").Value = "  +8.87%  "

# Row 20
$ws.Range("D20").Value = "525.94"
$ws.Range("E20").Value = "  +8.29%  "

# Row 21
$ws.Range("D21").Value = "16.92"
$ws.Range("E21").Value = "  +2.56%  "

# Row 22
$ws.Range("E22").Value = "  +23.72%  "

# Row 23
$ws.Range("D23").Value = "0.750"
$ws.Range("E23").Value = "  +8.92%  "

# Row 24
$ws.Range("D24").Value = "88.32"
$ws.Range("E24").Value = "  +5.66%  "

# Row 25
$ws.Range("D25").Value = "2.50"
$ws.Range("E25").Value = "  +11.20%  "

# Row 26
$ws.Range("D26").Value = "13.57"
$ws.Range("E26").Value = "  +7.05%  "

# Row 27
$ws.Range("D27").Value = "10.99"
$ws.Range("E27").Value = "  +6.65%  "

# Row 28
$ws.Range("E28").Value = "  -0.10%  "

# Row 29
$ws.Range("D29").Value = "0.0000123"
$ws.Range("E29").Value = "  +30.83%  "

# Row 30
$ws.Range("D30").Value = "2.54"
$ws.Range("E30").Value = "  +9.64%  "

# Row 31
$ws.Range("D31").Value = "2.89"
$ws.Range("E31").Value = "  +11.43%  "

# Row 32
$ws.Range("D32").Value = "7.96"
$ws.Range("E32").Value = "  +0.82%  "

# Row 33
$ws.Range("D33").Value = "32.22"
$ws.Range("E33").Value = "  +14.72%  "

# Row 34
$ws.Range("E34").Value = "  +3.47%  "

# Row 35
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.10%  "

# Row 36
$ws.Range("D36").Value = "6.20"
$ws.Range("E36").Value = "  +11.27%  "

# Row 37
$ws.Range("E37").Value = "  +10.42%  "

# Row 38
$ws.Range("D38").Value = "0.343"
$ws.Range("E38").Value = "  +9.43%  "

# Row 39
$ws.Range("D39").Value = "2.18"
$ws.Range("E39").Value = "  +8.58%  "

# Row 40
$ws.Range("E40").Value = "  +8.45%  "

# Row 41
$ws.Range("D41").Value = "51.60"
$ws.Range("E41").Value = "  +5.01%  "

# Row 42
$ws.Range("D42").Value = "8.91"
$ws.Range("E42").Value = "  +7.84%  "

# Row 43
$ws.Range("D43").Value = "44.56"
$ws.Range("E43").Value = "  -5.81%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "3.141.82"
$ws.Range("E44").Value = "  +12.81%  "

# Row 45
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "427.32"
$ws.Range("E45").Value = "  +15.49%  "

# Row 46
$ws.Range("D46").Value = "2.78"
$ws.Range("E46").Value = "  +3.04%  "

# Row 47
$ws.Range("E47").Value = "  +7.24%  "

# Row 48
$ws.Range("D48").Value = "27.94"
$ws.Range("E48").Value = "  +5.76%  "

# Row 49
$ws.Range("D49").Value = "138.81"
$ws.Range("E49").Value = "  +2.21%  "

# Row 50
$ws.Range("D50").Value = "2.55"
$ws.Range("E50").Value = "  +11.16%  "

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
